# Update leve profit calculations across multiple sheets (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 3750
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3750
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 11250
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -11488

$ws.Range("H58").Value = 2595.8
$ws.Range("I58").Value = 220
$ws.Range("J58").Value = 4674.625
$ws.Range("K58").Value = 660
$ws.Range("L58").Value = 14023.875
$ws.Range("M58").Value = -510
$ws.Range("N58").Value = -14323.875

$ws.Range("H60").Value = 3750
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 3750
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11250
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = -12218

$ws.Range("H108").Value = 39866
$ws.Range("J108").Value = 39866
$ws.Range("L108").Value = 39866
$ws.Range("N108").Value = -47546

$ws.Range("H129").Value = 1314.4423
$ws.Range("I129").Value = 697.7143
$ws.Range("J129").Value = 1410.3778
$ws.Range("K129").Value = 2093.1429
$ws.Range("L129").Value = 4231.1334
$ws.Range("M129").Value = 2906.8571
$ws.Range("N129").Value = -14231.1334

$ws.Range("H141").Value = 1727.0244
$ws.Range("I141").Value = 1212.4242
$ws.Range("J141").Value = 3849.75
$ws.Range("K141").Value = 3637.2726
$ws.Range("L141").Value = 11549.25
$ws.Range("M141").Value = 1542.7274
$ws.Range("N141").Value = -21909.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""

$ws.Range("H95").Value = 38543
$ws.Range("J95").Value = 38543
$ws.Range("L95").Value = 38543
$ws.Range("N95").Value = -44035

$ws.Range("H101").Value = 44426
$ws.Range("J101").Value = 44426
$ws.Range("L101").Value = 44426
$ws.Range("N101").Value = -50916

$ws.Range("H113").Value = 53266.668
$ws.Range("J113").Value = 53266.668
$ws.Range("L113").Value = 53266.668
$ws.Range("N113").Value = -61944.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3367.3333
$ws.Range("I5").Value = 480.8
$ws.Range("J5").Value = 17800
$ws.Range("K5").Value = 480.8
$ws.Range("L5").Value = 17800
$ws.Range("M5").Value = -367.8
$ws.Range("N5").Value = -18026

$ws.Range("H29").Value = 5333.3335
$ws.Range("I29").Value = 3000
$ws.Range("J29").Value = 6500
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 6500
$ws.Range("M29").Value = -2711
$ws.Range("N29").Value = -7078

$ws.Range("H107").Value = 3037.4443
$ws.Range("J107").Value = 4317.857
$ws.Range("L107").Value = 4317.857
$ws.Range("N107").Value = -8157.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1536
$ws.Range("I10").Value = 1536
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1536
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1397
$ws.Range("N10").Value = ""

$ws.Range("H17").Value = 27500
$ws.Range("I17").Value = 20000
$ws.Range("J17").Value = 35000
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 35000
$ws.Range("M17").Value = -19826
$ws.Range("N17").Value = -35348

$ws.Range("H22").Value = 422.75
$ws.Range("I22").Value = 397.42856
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 397.42856
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -47.42856
$ws.Range("N22").Value = -1300

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""

$ws.Range("H57").Value = 8000
$ws.Range("J57").Value = 8000
$ws.Range("L57").Value = 8000
$ws.Range("N57").Value = -9120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1111697.1
$ws.Range("I5").Value = 670.2222
$ws.Range("J5").Value = 4444777.5
$ws.Range("K5").Value = 2010.6666
$ws.Range("L5").Value = 13334332.5
$ws.Range("M5").Value = -1898.6666
$ws.Range("N5").Value = -13334556.5

$ws.Range("H64").Value = 6655.3335
$ws.Range("J64").Value = 10239.846
$ws.Range("L64").Value = 30719.538
$ws.Range("N64").Value = -31259.538

$ws.Range("H67").Value = 6655.3335
$ws.Range("J67").Value = 10239.846
$ws.Range("L67").Value = 30719.538
$ws.Range("N67").Value = -32591.538

$ws.Range("H88").Value = 4130.4443
$ws.Range("J88").Value = 4130.4443
$ws.Range("L88").Value = 12391.3329
$ws.Range("N88").Value = -13247.3329

$ws.Range("H91").Value = 4130.4443
$ws.Range("J91").Value = 4130.4443
$ws.Range("L91").Value = 12391.3329
$ws.Range("N91").Value = -15355.3329

$ws.Range("H113").Value = 622.60974
$ws.Range("I113").Value = 558.53845
$ws.Range("J113").Value = 733.6667
$ws.Range("K113").Value = 1675.61535
$ws.Range("L113").Value = 2201.0001
$ws.Range("M113").Value = 494.38465
$ws.Range("N113").Value = -6541.0001

$ws.Range("H131").Value = 893.8101
$ws.Range("I131").Value = 460
$ws.Range("J131").Value = 923.12164
$ws.Range("K131").Value = 1380
$ws.Range("L131").Value = 2769.36492
$ws.Range("M131").Value = 3660
$ws.Range("N131").Value = -12849.36492

$ws.Range("H135").Value = 1111697.1
$ws.Range("I135").Value = 670.2222
$ws.Range("J135").Value = 4444777.5
$ws.Range("K135").Value = 6031.999800000001
$ws.Range("L135").Value = 40002997.5
$ws.Range("M135").Value = -3496.999800000001
$ws.Range("N135").Value = -40008067.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 671.7692
$ws.Range("I107").Value = 540.1053000000001
$ws.Range("J107").Value = 1029.1428
$ws.Range("K107").Value = 540.1053000000001
$ws.Range("L107").Value = 1029.1428
$ws.Range("M107").Value = 1379.8947
$ws.Range("N107").Value = -4869.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 22752
$ws.Range("I17").Value = 24669.334
$ws.Range("J17").Value = 17000
$ws.Range("K17").Value = 24669.334
$ws.Range("L17").Value = 17000
$ws.Range("M17").Value = -24497.334
$ws.Range("N17").Value = -17344

$ws.Range("H63").Value = 16613
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248

$ws.Range("H66").Value = 16613
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
